$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: [008-0001] paragraph -- merge "[008-0001]" + (bookmark) + "Add
# product to checkout list" into a single run with no bookmark in between.
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$p1 = $null
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.StartsWith("[008-0001]")) {
        $p1 = $p
        break
    }
}

if ($p1 -ne $null) {
    # Setting Range.Text to the exact same visible text is a no-op for the
    # engine (it short-circuits on equal text), so first stamp a distinct
    # placeholder, then set the real text -- this forces a genuine
    # delete+insert that merges the run and drops the old bookmark that
    # was sitting inside the replaced range.
    $full = $p1.Range
    $body = $d.Range($full.Start, $full.End - 1)
    $body.Text = "TEMP_PLACEHOLDER_0001"

    $full2 = $p1.Range
    $body2 = $d.Range($full2.Start, $full2.End - 1)
    $body2.Text = "[008-0001]Add product to checkout list"
}

# ---------------------------------------------------------------------------
# Step 2: [008-0013] paragraph -- insert "invalid " before "numeric" and move
# the _GoBack bookmark to sit right before "numeric ".
# ---------------------------------------------------------------------------
$p2 = $null
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.StartsWith("[008-0013]")) {
        $p2 = $p
        break
    }
}

if ($p2 -ne $null) {
    $rng = $p2.Range
    $rng.Find.Execute("with numeric", $false, $false, $false, $false, $false, $true, 1, $false, "with invalid numeric", 2) | Out-Null

    $rng2 = $p2.Range
    $rng2.Find.Execute("numeric", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

    $bmRange = $d.Range($rng2.Start, $rng2.Start)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Output "edit complete"
